$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 23811238
$ws.Range("J112").Value = 1788.7179
$ws.Range("L112").Value = 5366.153700000001
$ws.Range("N112").Value = -7582.153700000001
$ws.Range("H132").Value = 35721356
$ws.Range("I132").Value = 40006720
$ws.Range("K132").Value = 120020160
$ws.Range("M132").Value = -120017630

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 660
$ws.Range("I2").Value = 668
$ws.Range("J2").Value = 630
$ws.Range("K2").Value = 668
$ws.Range("L2").Value = 630
$ws.Range("M2").Value = -555
$ws.Range("N2").Value = -856
$ws.Range("H61").Value = 1182.8235
$ws.Range("I61").Value = 927.6429000000001
$ws.Range("K61").Value = 927.6429000000001
$ws.Range("M61").Value = -715.6429000000001
$ws.Range("H74").Value = 3168.805
$ws.Range("I74").Value = 3176.9092
$ws.Range("J74").Value = 3135.375
$ws.Range("K74").Value = 3176.9092
$ws.Range("L74").Value = 3135.375
$ws.Range("M74").Value = -2302.9092
$ws.Range("N74").Value = -4883.375
$ws.Range("H77").Value = 3168.805
$ws.Range("I77").Value = 3176.9092
$ws.Range("J77").Value = 3135.375
$ws.Range("K77").Value = 15884.546
$ws.Range("L77").Value = 15676.875
$ws.Range("M77").Value = -11516.546
$ws.Range("N77").Value = -24412.875
$ws.Range("H116").Value = 660
$ws.Range("I116").Value = 668
$ws.Range("J116").Value = 630
$ws.Range("K116").Value = 668
$ws.Range("L116").Value = 630
$ws.Range("M116").Value = 1626
$ws.Range("N116").Value = -5218
$ws.Range("H136").Value = 1182.8235
$ws.Range("I136").Value = 927.6429000000001
$ws.Range("K136").Value = 2782.9287
$ws.Range("M136").Value = -232.9287000000004

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 660
$ws.Range("I3").Value = 668
$ws.Range("J3").Value = 630
$ws.Range("K3").Value = 668
$ws.Range("L3").Value = 630
$ws.Range("M3").Value = -554
$ws.Range("N3").Value = -858
$ws.Range("H107").Value = 1559.9
$ws.Range("I107").Value = 1324.875
$ws.Range("K107").Value = 1324.875
$ws.Range("M107").Value = 595.125
$ws.Range("H134").Value = 1783.1951
$ws.Range("I134").Value = 1111.6167
$ws.Range("J134").Value = 3614.7727
$ws.Range("K134").Value = 3334.8501
$ws.Range("L134").Value = 10844.3181
$ws.Range("M134").Value = -799.8501000000001
$ws.Range("N134").Value = -15914.3181
$ws.Range("H141").Value = 31000
$ws.Range("J141").Value = 31000
$ws.Range("L141").Value = 31000
$ws.Range("N141").Value = -41360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5292808.5
$ws.Range("I16").Value = 9260695
$ws.Range("J16").Value = 2293
$ws.Range("K16").Value = 9260695
$ws.Range("L16").Value = 2293
$ws.Range("M16").Value = -9260408
$ws.Range("N16").Value = -2867
$ws.Range("H31").Value = 7248683
$ws.Range("I31").Value = 1245.6
$ws.Range("J31").Value = 26320886
$ws.Range("K31").Value = 1245.6
$ws.Range("L31").Value = 26320886
$ws.Range("M31").Value = -950.5999999999999
$ws.Range("N31").Value = -26321476
$ws.Range("H34").Value = 7248683
$ws.Range("I34").Value = 1245.6
$ws.Range("J34").Value = 26320886
$ws.Range("K34").Value = 1245.6
$ws.Range("L34").Value = 26320886
$ws.Range("M34").Value = -1043.6
$ws.Range("N34").Value = -26321290
$ws.Range("H58").Value = 1669.2106
$ws.Range("I58").Value = 1585.6617
$ws.Range("J58").Value = 1879.6296
$ws.Range("K58").Value = 1585.6617
$ws.Range("L58").Value = 1879.6296
$ws.Range("M58").Value = -1382.6617
$ws.Range("N58").Value = -2285.6296
$ws.Range("H113").Value = 5292808.5
$ws.Range("I113").Value = 9260695
$ws.Range("J113").Value = 2293
$ws.Range("K113").Value = 9260695
$ws.Range("L113").Value = 2293
$ws.Range("M113").Value = -9258525
$ws.Range("N113").Value = -6633
$ws.Range("H132").Value = 3372.7307
$ws.Range("I132").Value = 3180.5625
$ws.Range("J132").Value = 3680.2
$ws.Range("K132").Value = 9541.6875
$ws.Range("L132").Value = 11040.6
$ws.Range("M132").Value = -7011.6875
$ws.Range("N132").Value = -16100.6
$ws.Range("H134").Value = 4261.488
$ws.Range("I134").Value = 5921.05
$ws.Range("K134").Value = 17763.15
$ws.Range("M134").Value = -15228.15
$ws.Range("H136").Value = 1669.2106
$ws.Range("I136").Value = 1585.6617
$ws.Range("J136").Value = 1879.6296
$ws.Range("K136").Value = 4756.9851
$ws.Range("L136").Value = 5638.8888
$ws.Range("M136").Value = -2206.9851
$ws.Range("N136").Value = -10738.8888

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 3080
$ws.Range("I123").Value = 4333.3335
$ws.Range("J123").Value = 1200
$ws.Range("K123").Value = 13000.0005
$ws.Range("L123").Value = 3600
$ws.Range("M123").Value = -10550.0005
$ws.Range("N123").Value = -8500
$ws.Range("H129").Value = 3101.5293
$ws.Range("I129").Value = 4105
$ws.Range("J129").Value = 2554.182
$ws.Range("K129").Value = 12315
$ws.Range("L129").Value = 7662.545999999999
$ws.Range("M129").Value = -7315
$ws.Range("N129").Value = -17662.546
$ws.Range("H130").Value = 3120
$ws.Range("I130").Value = 1695
$ws.Range("J130").Value = 4260
$ws.Range("K130").Value = 5085
$ws.Range("L130").Value = 12780
$ws.Range("M130").Value = -65
$ws.Range("N130").Value = -22820
$ws.Range("H131").Value = 824.8570999999999
$ws.Range("J131").Value = 874.6716300000001
$ws.Range("L131").Value = 2624.01489
$ws.Range("N131").Value = -12704.01489
$ws.Range("H133").Value = 3369.5
$ws.Range("J133").Value = 2635.3845
$ws.Range("L133").Value = 7906.1535
$ws.Range("N133").Value = -18026.1535
$ws.Range("H134").Value = 4165.4243
$ws.Range("I134").Value = 3570.5
$ws.Range("K134").Value = 10711.5
$ws.Range("M134").Value = -5641.5
$ws.Range("H136").Value = 3344.0588
$ws.Range("I136").Value = 3089.2144
$ws.Range("J136").Value = 4533.3335
$ws.Range("K136").Value = 9267.643199999999
$ws.Range("L136").Value = 13600.0005
$ws.Range("M136").Value = -4167.643199999999
$ws.Range("N136").Value = -23800.0005
$ws.Range("H137").Value = 2450.3635
$ws.Range("J137").Value = 2978.125
$ws.Range("L137").Value = 8934.375
$ws.Range("N137").Value = -19134.375
$ws.Range("H138").Value = 2760.8333
$ws.Range("H139").Value = 1562.2354
$ws.Range("I139").Value = 1039.8572
$ws.Range("K139").Value = 3119.5716
$ws.Range("M139").Value = 2020.4284
$ws.Range("H140").Value = 2778.7778
$ws.Range("I140").Value = 3180.5
$ws.Range("K140").Value = 9541.5
$ws.Range("M140").Value = -4361.5
$ws.Range("H141").Value = 10872.223
$ws.Range("I141").Value = 12321.429
$ws.Range("K141").Value = 36964.287
$ws.Range("M141").Value = -31784.287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2891.62
$ws.Range("I126").Value = 2891.62
$ws.Range("K126").Value = 8674.860000000001
$ws.Range("M126").Value = -6204.860000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2565.5
$ws.Range("I100").Value = 2131
$ws.Range("K100").Value = 2131
$ws.Range("M100").Value = -1590
$ws.Range("H132").Value = 3888.7742
$ws.Range("I132").Value = 1327.775
$ws.Range("J132").Value = 8545.137000000001
$ws.Range("K132").Value = 3983.325
$ws.Range("L132").Value = 25635.411
$ws.Range("M132").Value = -1453.325
$ws.Range("N132").Value = -30695.411

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 238251.67
$ws.Range("I126").Value = 1205.2
$ws.Range("J126").Value = 534559.75
$ws.Range("K126").Value = 3615.6
$ws.Range("L126").Value = 1603679.25
$ws.Range("M126").Value = -1145.6
$ws.Range("N126").Value = -1608619.25
$ws.Range("H136").Value = 1358.25
$ws.Range("I136").Value = 789.1667
$ws.Range("K136").Value = 2367.5001
$ws.Range("M136").Value = 182.4998999999998
